$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Update column A labels for rows 5-10: rename the "Audiodescriptor1..6" entries
# to the correct descriptor names (bug fix after using convert.py again).
$ws.Range("A5").Value = "Rms"
$ws.Range("A6").Value = "Pitch"
$ws.Range("A7").Value = "Centroid"
$ws.Range("A8").Value = "Flatness"
$ws.Range("A9").Value = "TransDens"
$ws.Range("A10").Value = "Spread"

# Update the view: move the active selection to A11 and scroll so that
# row 5 is the top-left visible row.
$ws.Range("A11").Select()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
